$d = $word.ActiveDocument

# 1. Thesis title: collapse the space before "Mono No Aware..." into ": "
$d.Content.Find.Execute(
    "”" + " " + "Mono No Aware and the Tanuki in",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "”" + ": Mono No Aware and the Tanuki in",
    2)

# 2. "Guile" -> "Guile Scheme" in the languages list
$d.Content.Find.Execute(
    "Guile",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Guile Scheme",
    2)

# 3. "Exclamat!ion" -> "Exclamat!on" (hyperlink text)
$d.Content.Find.Execute(
    "Exclamat!ion",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Exclamat!on",
    2)
